# Add a new ticket row (row 22) to the tickets worksheet, matching the
# existing "Laptop" / "Hardware" ticket rows above it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "TCKT-1021"
$ws.Range("B22").Value = "Laptop battery issue"
$ws.Range("C22").Value = "Hardware"
$ws.Range("D22").Value = "Laptop"
$ws.Range("E22").Value = "Check charger, if charger is working fine then replace with new battery"

# Leave the newly entered cell selected, as in the saved workbook.
$null = $ws.Range("B22").Select()
